$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mgz1_antibiotic_11.9")

# Fix casing of the dCasRx strain label (was "dCASRx") on the wells that use it
$ws.Range("B5").Value = "dCasRx"
$ws.Range("B6").Value = "dCasRx"
$ws.Range("B7").Value = "dCasRx"
$ws.Range("B11").Value = "dCasRx"
$ws.Range("B12").Value = "dCasRx"
$ws.Range("B13").Value = "dCasRx"
$ws.Range("B17").Value = "dCasRx"
$ws.Range("B18").Value = "dCasRx"
$ws.Range("B19").Value = "dCasRx"

# Correct the compound_1 concentration unit from mM to uM for these wells
$ws.Range("H14:H19").Value = "uM"

# Move the active selection to H19, matching where the author left off
$ws.Range("H19").Select() | Out-Null
